# Add a new "Feature 2 Added" slide, modeled on the existing
# "Feature Added" slide (slide 2): duplicate it so the new slide
# inherits the same shape/text formatting, then update its text.

$p = $ppt.ActivePresentation

$sourceSlide = $p.Slides.Item(2)
$newSlideRange = $sourceSlide.Duplicate()
$newSlide = $newSlideRange.Item(1)

$textBox = $newSlide.Shapes.Item(1)
$textRange = $textBox.TextFrame.TextRange

$textRange.Paragraphs(1).Runs(1).Text = "Feature 2 Added"
$textRange.Paragraphs(3).Runs(1).Text = "04/04/2020"
